$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 238.66667
$ws.Range("I12").Value = 48.5
$ws.Range("J12").Value = 333.75
$ws.Range("K12").Value = 48.5
$ws.Range("L12").Value = 333.75
$ws.Range("M12").Value = 121.5
$ws.Range("N12").Value = -673.75
# Row 33
$ws.Range("H33").Value = 263.53333
$ws.Range("I33").Value = 132.81818
$ws.Range("K33").Value = 132.81818
$ws.Range("M33").Value = 96.18181999999999
# Row 38
$ws.Range("H38").Value = 33.333332
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
# Row 62
$ws.Range("H62").Value = 3758.6428
$ws.Range("I62").Value = 3758.6428
$ws.Range("K62").Value = 3758.6428
$ws.Range("M62").Value = -3134.6428
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
# Row 65
$ws.Range("H65").Value = 3758.6428
$ws.Range("I65").Value = 3758.6428
$ws.Range("K65").Value = 18793.214
$ws.Range("M65").Value = -15673.214
# Row 67
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
# Row 100
$ws.Range("H100").Value = 3183.3333
$ws.Range("I100").Value = 3333.3333
$ws.Range("K100").Value = 3333.3333
$ws.Range("M100").Value = -2792.3333
# Row 103
$ws.Range("H103").Value = 1881.3334
$ws.Range("I103").Value = 1897.6
$ws.Range("J103").Value = 1800
$ws.Range("K103").Value = 5692.799999999999
$ws.Range("L103").Value = 5400
$ws.Range("M103").Value = -5106.799999999999
$ws.Range("N103").Value = -6572
# Row 132
$ws.Range("H132").Value = 2262.4736
$ws.Range("I132").Value = 2262.4736
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6787.4208
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4257.4208
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2939.25
$ws.Range("I32").Value = 2722.0715
$ws.Range("K32").Value = 2722.0715
$ws.Range("M32").Value = -2435.0715
# Row 74
$ws.Range("H74").Value = 1146
$ws.Range("I74").Value = 1160.8334
$ws.Range("K74").Value = 1160.8334
$ws.Range("M74").Value = -286.8334
# Row 77
$ws.Range("H77").Value = 1146
$ws.Range("I77").Value = 1160.8334
$ws.Range("K77").Value = 5804.166999999999
$ws.Range("M77").Value = -1436.166999999999
# Row 88
$ws.Range("H88").Value = 3809.0908
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3809.0908
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3809.0908
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -4621.0908
# Row 91
$ws.Range("H91").Value = 3809.0908
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3809.0908
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3809.0908
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -6617.0908
# Row 122
$ws.Range("H122").Value = 2867.3333
$ws.Range("I122").Value = 2940.9
$ws.Range("J122").Value = 2499.5
$ws.Range("K122").Value = 8822.700000000001
$ws.Range("L122").Value = 7498.5
$ws.Range("M122").Value = -6372.700000000001
$ws.Range("N122").Value = -12398.5
# Row 132
$ws.Range("H132").Value = 2302.7778
$ws.Range("I132").Value = 2302.7778
$ws.Range("K132").Value = 6908.3334
$ws.Range("M132").Value = -4378.3334

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 3679.889
$ws.Range("I105").Value = 3515
$ws.Range("K105").Value = 3515
$ws.Range("M105").Value = -1768

$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Range("H134").Value = 5999.5
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()
# Row 138
$ws.Range("H138").Value = 49499.5
$ws.Range("J138").Value = 49499.5
$ws.Range("L138").Value = 49499.5
$ws.Range("N138").Value = -59779.5

$ws = $wb.Worksheets.Item("CUL")
# Row 26
$ws.Range("H26").Value = 176016.17
$ws.Range("I26").Value = 251450
$ws.Range("J26").Value = 25148.5
$ws.Range("K26").Value = 754350
$ws.Range("L26").Value = 75445.5
$ws.Range("M26").Value = -754062
$ws.Range("N26").Value = -76021.5
# Row 81
$ws.Range("H81").Value = 4166.3335
$ws.Range("I81").Value = 4499
$ws.Range("J81").Value = 4000
$ws.Range("K81").Value = 13497
$ws.Range("L81").Value = 12000
$ws.Range("M81").Value = -12374
$ws.Range("N81").Value = -14246
# Row 84
$ws.Range("H84").Value = 4166.3335
$ws.Range("I84").Value = 4499
$ws.Range("J84").Value = 4000
$ws.Range("K84").Value = 40491
$ws.Range("L84").Value = 36000
$ws.Range("M84").Value = -34875
$ws.Range("N84").Value = -47232
# Row 92
$ws.Range("H92").Value = 636.2
$ws.Range("J92").Value = 781.3333
$ws.Range("L92").Value = 2343.9999
$ws.Range("N92").Value = -4839.9999
# Row 117
$ws.Range("H117").Value = 545.25
$ws.Range("I117").Value = 499.33334
$ws.Range("K117").Value = 1498.00002
$ws.Range("M117").Value = 1943.99998

$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 14166.667
# Row 57
$ws.Range("H57").Value = 23657.889
$ws.Range("J57").Value = 25115.125
$ws.Range("L57").Value = 25115.125
$ws.Range("N57").Value = -26755.125
# Row 126
$ws.Range("H126").Value = 12772.846
$ws.Range("I126").Value = 13281
$ws.Range("J126").Value = 11959.8
$ws.Range("K126").Value = 39843
$ws.Range("L126").Value = 35879.39999999999
$ws.Range("M126").Value = -37373
$ws.Range("N126").Value = -40819.39999999999
# Row 132
$ws.Range("H132").Value = 4189.3887
$ws.Range("I132").Value = 4088.0625
$ws.Range("K132").Value = 12264.1875
$ws.Range("M132").Value = -9734.1875
# Row 134
$ws.Range("H134").Value = 149998.5
$ws.Range("J134").Value = 149998.5
$ws.Range("L134").Value = 449995.5
$ws.Range("N134").Value = -455065.5
# Row 136
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 5
$ws.Range("H5").Value = 1325
$ws.Range("I5").Value = 1325
$ws.Range("K5").Value = 1325
$ws.Range("M5").Value = -1212
# Row 22
$ws.Range("H22").Value = 999
$ws.Range("I22").Value = 999
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 999
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -704
$ws.Range("N22").ClearContents()
# Row 27
$ws.Range("H27").Value = 999
$ws.Range("I27").Value = 999
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 999
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -892
$ws.Range("N27").ClearContents()
# Row 40
$ws.Range("H40").Value = 1741
$ws.Range("J40").Value = 1697
$ws.Range("L40").Value = 1697
$ws.Range("N40").Value = -1969
# Row 46
$ws.Range("H46").Value = 2641.88
$ws.Range("I46").Value = 2176.353
$ws.Range("J46").Value = 3631.125
$ws.Range("K46").Value = 2176.353
$ws.Range("L46").Value = 3631.125
$ws.Range("M46").Value = -1988.353
$ws.Range("N46").Value = -4007.125
# Row 122
$ws.Range("H122").Value = 3402
$ws.Range("I122").Value = 3402
$ws.Range("K122").Value = 10206
$ws.Range("M122").Value = -7756
# Row 132
$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2659.7778
$ws.Range("I122").Value = 2499.7856
$ws.Range("K122").Value = 7499.3568
$ws.Range("M122").Value = -5049.3568
# Row 132
$ws.Range("H132").Value = 1622.1111
$ws.Range("I132").Value = 1524.75
$ws.Range("K132").Value = 4574.25
$ws.Range("M132").Value = -2044.25
